$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K11").Value = -0.0703
$ws.Range("L11").Value = -0.0731
$ws.Range("M11").Value = -0.055
$ws.Range("N11").Value = -0.0487
$ws.Range("O11").Value = -0.0394
$ws.Range("P11").Value = -0.0275
$ws.Range("Q11").Value = -0.0246

$ws.Range("O15").Value = 0.1478
$ws.Range("P15").Value = 0.1643
$ws.Range("Q15").Value = 0.2998

$ws.Range("J29").Value = -2.0918
$ws.Range("K29").Value = -2.087
$ws.Range("L29").Value = -2.9915
$ws.Range("M29").Value = -1.2599
$ws.Range("N29").Value = -0.9552
$ws.Range("O29").Value = -0.4207
$ws.Range("P29").Value = -0.5797
$ws.Range("Q29").Value = -0.4185

$ws.Range("J43").Value = -0.0172
$ws.Range("K43").Value = -0.0195
$ws.Range("L43").Value = -0.0206
$ws.Range("M43").Value = -0.0122
$ws.Range("N43").Value = -0.0104
$ws.Range("O43").Value = -0.0067
$ws.Range("P43").Value = -0.0037
$ws.Range("Q43").Value = -0.002

$ws.Range("J47").Value = -0.1418
$ws.Range("K47").Value = -0.0963
$ws.Range("L47").Value = -0.0672
$ws.Range("M47").Value = -0.062
$ws.Range("N47").Value = -0.0451
$ws.Range("O47").Value = -0.0348
$ws.Range("P47").Value = -0.0272
$ws.Range("Q47").Value = -0.0009

$ws.Range("K65").Value = -0.0268
$ws.Range("L65").Value = -0.0265
$ws.Range("M65").Value = -0.0262
$ws.Range("N65").Value = -0.0258
$ws.Range("O65").Value = -0.0254
$ws.Range("P65").Value = -0.025
$ws.Range("Q65").Value = -0.0247

$ws.Range("O69").Value = 0.0037
$ws.Range("P69").Value = 0.011
$ws.Range("Q69").Value = 0.0054

$ws.Range("J83").Value = 0.2094
$ws.Range("K83").Value = 0.1077
$ws.Range("L83").Value = 0.0725
$ws.Range("M83").Value = 0.0262
$ws.Range("N83").Value = -0.073
$ws.Range("O83").Value = -0.0583
$ws.Range("P83").Value = -0.0495
$ws.Range("Q83").Value = -0.0496

$ws.Range("J97").Value = -0.0052
$ws.Range("K97").Value = -0.005
$ws.Range("L97").Value = -0.0049
$ws.Range("M97").Value = -0.0047
$ws.Range("N97").Value = -0.0046
$ws.Range("O97").Value = -0.0045
$ws.Range("P97").Value = -0.0044
$ws.Range("Q97").Value = -0.0043

$ws.Range("J101").Value = -0.01
$ws.Range("K101").Value = -0.005
$ws.Range("L101").Value = -0.0049
$ws.Range("M101").Value = -0.0048
$ws.Range("N101").Value = -0.0048
$ws.Range("O101").Value = -0.0047
$ws.Range("P101").Value = -0.0046
$ws.Range("Q101").Value = 0
